$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values (G2, H2, G3, H3, G5, H5); I and J columns are
# formulas that recalculate automatically.
$ws.Range("G2").Value = 0.00621917999999999
$ws.Range("H2").Value = 1.961046

$ws.Range("G3").Value = 0.000412851
$ws.Range("H3").Value = 0.10091772

$ws.Range("G5").Value = 0.018693
$ws.Range("H5").Value = 0.16778596875

# Update the active selection on the sheet view.
$ws.Range("E31").Select()
